# Natmi following Dr Hou advice
# Adds a new "FAPs" cluster/string and fills rows 2-7 of Sheet1 with the
# recomputed ligand-receptor pair statistics for every sCs/FAPs combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "sCs"
$ws.Cells.Item(2, 2).Value2 = "Edn3"
$ws.Cells.Item(2, 3).Value2 = "Ednrb"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 7).Value2 = 0.356007
$ws.Cells.Item(2, 8).Value2 = 1.068021
$ws.Cells.Item(2, 9).Value2 = 0.08840121110588733
$ws.Cells.Item(2, 10).Value2 = 0.08840121110588735
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 21.22137533333333
$ws.Cells.Item(2, 14).Value2 = 63.664126
$ws.Cells.Item(2, 15).Value2 = 0.2719819326156953
$ws.Cells.Item(2, 16).Value2 = 0.2719819326156953
$ws.Cells.Item(2, 17).Value2 = 7.554958168294001
$ws.Cells.Item(2, 18).Value2 = 67.99462351464601
$ws.Cells.Item(2, 19).Value2 = 0.02404353224214731
$ws.Cells.Item(2, 20).Value2 = 0.02404353224214731

# Row 3
$ws.Cells.Item(3, 1).Value2 = "sCs"
$ws.Cells.Item(3, 2).Value2 = "Edn3"
$ws.Cells.Item(3, 3).Value2 = "Ednrb"
$ws.Cells.Item(3, 4).Value2 = "sCs"
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 7).Value2 = 0.356007
$ws.Cells.Item(3, 8).Value2 = 1.068021
$ws.Cells.Item(3, 9).Value2 = 0.08840121110588733
$ws.Cells.Item(3, 10).Value2 = 0.08840121110588735
$ws.Cells.Item(3, 11).Value2 = 1
$ws.Cells.Item(3, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 13).Value2 = 0.01207966666666667
$ws.Cells.Item(3, 14).Value2 = 0.036239
$ws.Cells.Item(3, 15).Value2 = 0.0001548180093772148
$ws.Cells.Item(3, 16).Value2 = 0.0001548180093772148
$ws.Cells.Item(3, 17).Value2 = 0.004300445891
$ws.Cells.Item(3, 18).Value2 = 0.03870401301900001
$ws.Cells.Item(3, 19).Value2 = 0.00001368609952994841
$ws.Cells.Item(3, 20).Value2 = 0.00001368609952994841

# Row 4
$ws.Cells.Item(4, 1).Value2 = "sCs"
$ws.Cells.Item(4, 2).Value2 = "Edn3"
$ws.Cells.Item(4, 3).Value2 = "Ednrb"
$ws.Cells.Item(4, 4).Value2 = "FAPs"
$ws.Cells.Item(4, 5).Value2 = 1
$ws.Cells.Item(4, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 7).Value2 = 0.356007
$ws.Cells.Item(4, 8).Value2 = 1.068021
$ws.Cells.Item(4, 9).Value2 = 0.08840121110588733
$ws.Cells.Item(4, 10).Value2 = 0.08840121110588735
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 56.79149
$ws.Cells.Item(4, 14).Value2 = 170.37447
$ws.Cells.Item(4, 15).Value2 = 0.7278632493749275
$ws.Cells.Item(4, 16).Value2 = 0.7278632493749274
$ws.Cells.Item(4, 17).Value2 = 20.21816798043
$ws.Cells.Item(4, 18).Value2 = 181.96351182387
$ws.Cells.Item(4, 19).Value2 = 0.06434399276421009
$ws.Cells.Item(4, 20).Value2 = 0.06434399276421009

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Edn3"
$ws.Cells.Item(5, 3).Value2 = "Ednrb"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 3.671166333333333
$ws.Cells.Item(5, 8).Value2 = 11.013499
$ws.Cells.Item(5, 9).Value2 = 0.9115987888941126
$ws.Cells.Item(5, 10).Value2 = 0.9115987888941127
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 21.22137533333333
$ws.Cells.Item(5, 14).Value2 = 63.664126
$ws.Cells.Item(5, 15).Value2 = 0.2719819326156953
$ws.Cells.Item(5, 16).Value2 = 0.2719819326156953
$ws.Cells.Item(5, 17).Value2 = 77.90719867076378
$ws.Cells.Item(5, 18).Value2 = 701.164788036874
$ws.Cells.Item(5, 19).Value2 = 0.247938400373548
$ws.Cells.Item(5, 20).Value2 = 0.247938400373548

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Edn3"
$ws.Cells.Item(6, 3).Value2 = "Ednrb"
$ws.Cells.Item(6, 4).Value2 = "sCs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 3.671166333333333
$ws.Cells.Item(6, 8).Value2 = 11.013499
$ws.Cells.Item(6, 9).Value2 = 0.9115987888941126
$ws.Cells.Item(6, 10).Value2 = 0.9115987888941127
$ws.Cells.Item(6, 11).Value2 = 1
$ws.Cells.Item(6, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(6, 13).Value2 = 0.01207966666666667
$ws.Cells.Item(6, 14).Value2 = 0.036239
$ws.Cells.Item(6, 15).Value2 = 0.0001548180093772148
$ws.Cells.Item(6, 16).Value2 = 0.0001548180093772148
$ws.Cells.Item(6, 17).Value2 = 0.04434646558455556
$ws.Cells.Item(6, 18).Value2 = 0.399118190261
$ws.Cells.Item(6, 19).Value2 = 0.0001411319098472664
$ws.Cells.Item(6, 20).Value2 = 0.0001411319098472664

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Edn3"
$ws.Cells.Item(7, 3).Value2 = "Ednrb"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 3.671166333333333
$ws.Cells.Item(7, 8).Value2 = 11.013499
$ws.Cells.Item(7, 9).Value2 = 0.9115987888941126
$ws.Cells.Item(7, 10).Value2 = 0.9115987888941127
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 56.79149
$ws.Cells.Item(7, 14).Value2 = 170.37447
$ws.Cells.Item(7, 15).Value2 = 0.7278632493749275
$ws.Cells.Item(7, 16).Value2 = 0.7278632493749274
$ws.Cells.Item(7, 17).Value2 = 208.4910061078367
$ws.Cells.Item(7, 18).Value2 = 1876.41905497053
$ws.Cells.Item(7, 19).Value2 = 0.6635192566107174
$ws.Cells.Item(7, 20).Value2 = 0.6635192566107174

